# Add Jimenez 2019 pattern generator and evaluation functions:
# update the generated pattern / evaluation table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 ---
$ws.Range("B4").Value = "TAATATATAT"
$ws.Range("H4").Value = "'1"
$ws.Range("H4").ClearFormats()

# --- Row 5 ---
$ws.Range("B5").Value = "ATATTAATAT"
$ws.Range("H5").Value = "'5"
$ws.Range("H5").ClearFormats()

# --- Row 10 ---
$ws.Range("B10").Value = "ATATAATAAT"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("G10").Value = "'5"
$ws.Range("G10").ClearFormats()
$ws.Range("H10").Value = "'7"
$ws.Range("H10").ClearFormats()

# --- Row 11 ---
$ws.Range("B11").Value = "AATAATATAT"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("G11").Value = "'1"
$ws.Range("G11").ClearFormats()
$ws.Range("H11").Value = "'3"
$ws.Range("H11").ClearFormats()

# --- Row 12 ---
$ws.Range("B12").Value = "TATAATATAT"
$ws.Range("H12").Value = "1, 3"

# --- Row 13 ---
$ws.Range("B13").Value = "ATATTATAAT"
$ws.Range("H13").Value = "5, 7"

# --- Row 14 ---
$ws.Range("B14").Value = "TAATATTAAT"
$ws.Range("H14").Value = "1, 7"

# --- Row 15 ---
$ws.Range("B15").Value = "TAATTAATAT"
$ws.Range("H15").Value = "1, 5"

# --- Row 16 ---
$ws.Range("B16").Value = "ATTAATTAAT"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = ""
$ws.Range("H16").Value = "3, 7"

# --- Row 17 ---
$ws.Range("B17").Value = "TAATATTAAT"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 2
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = "1, 7"

# --- Row 18 ---
$ws.Range("B18").Value = "TAATAATAAT"
$ws.Range("G18").Value = "'5"
$ws.Range("G18").ClearFormats()
$ws.Range("H18").Value = "1, 7"

# --- Row 19 ---
$ws.Range("B19").Value = "AATAATTAAT"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2
$ws.Range("G19").Value = "'1"
$ws.Range("G19").ClearFormats()
$ws.Range("H19").Value = "3, 7"

# --- Row 20 ---
$ws.Range("B20").Value = "TAATTATAAT"
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = "1, 5, 7"
